$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Re-sort the table's data range (A3:T42) ascending by column A ("Id"),
# matching the updated sortCondition (was column H, now column A).
$lo = $ws.ListObjects.Item(1)
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("A3:A42"), 0, 1, 0, 0) | Out-Null
$lo.Sort.Header = 1
$lo.Sort.Apply()

# Update the saved selection on the active (bottom-right, frozen-pane) pane.
$ws.Activate()
$ws.Range("G10").Select()
